# PolicyBazar_UserData.xlsx edit:
#  - header row: "Name"/"Email" -> "Company"/"Premium"; "Phone Number" and
#    "Travel Country" headers cleared (fill recolored to white)
#  - data row 2 cleared entirely (name/email/phone/country removed)
#  - mailto: hyperlink on the old email cell removed
#  - selection moved to P12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the two headers that survive.
$ws.Range("A1").Value = "Company"
$ws.Range("B1").Value = "Premium"

# Drop the hyperlink that lived on B2 (mailto: link + its underlying text).
[void]$ws.Hyperlinks.Delete()

# Clear the now-unused header cells, but repaint them with a plain
# white fill (was the same orange/teal themed fill as A1:B1).
$ws.Range("C1:D1").ClearContents()
$ws.Range("C1:D1").Interior.ThemeColor = 2

# Clear every cell of the old sample data row.
$ws.Range("A2:D2").ClearContents()

# Match the author's final selection.
[void]$ws.Range("P12").Select()
